$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MinMaxDed")
$ws.Range("B1").Value = $ws.Range("B1").Value
